$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.043564110620407
$ws.Range("D2").Value = 1.057437346453313
$ws.Range("E2").Value = 1.052412955094316
$ws.Range("F2").Value = 1.065300448093593
$ws.Range("I2").Value = 1.049309076013587
$ws.Range("J2").Value = 1.048633882364249
$ws.Range("K2").Value = 1.060172318341023
$ws.Range("L2").Value = 1.055161757196498
$ws.Range("M2").Value = 1.068014059272366
$ws.Range("N2").Value = 1.050123062757194

$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.04443893259701
$ws.Range("D3").Value = 1.058058656770746
$ws.Range("E3").Value = 1.053163349265108
$ws.Range("F3").Value = 1.066067987594076
$ws.Range("I3").Value = 1.049536468800306
$ws.Range("J3").Value = 1.049155973309636
$ws.Range("K3").Value = 1.060607954033397
$ws.Range("L3").Value = 1.055725152340717
$ws.Range("M3").Value = 1.068597090606235
$ws.Range("N3").Value = 1.050645895131608

$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.045005659924496
$ws.Range("D4").Value = 1.058461066496904
$ws.Range("E4").Value = 1.053649811954028
$ws.Range("F4").Value = 1.06656545401872
$ws.Range("I4").Value = 1.049682583104031
$ws.Range("J4").Value = 1.049493796657062
$ws.Range("K4").Value = 1.060889525731639
$ws.Range("L4").Value = 1.056089946498114
$ws.Range("M4").Value = 1.068974502373305
$ws.Range("N4").Value = 1.050984198226936

$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.045244068299901
$ws.Range("D5").Value = 1.058630328756671
$ws.Range("E5").Value = 1.053854536295352
$ws.Range("F5").Value = 1.066774782543493
$ws.Range("I5").Value = 1.049743763547293
$ws.Range("J5").Value = 1.049635815260602
$ws.Range("K5").Value = 1.061007821999007
$ws.Range("L5").Value = 1.056243361949563
$ws.Range("M5").Value = 1.069133200967792
$ws.Range("N5").Value = 1.051126418513183

$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.045284107229902
$ws.Range("D6").Value = 1.058658753812133
$ws.Range("E6").Value = 1.053888922992698
$ws.Range("F6").Value = 1.06680994102194
$ws.Range("I6").Value = 1.049754021561969
$ws.Range("J6").Value = 1.049659660635196
$ws.Range("K6").Value = 1.061027679938518
$ws.Range("L6").Value = 1.056269124320563
$ws.Range("M6").Value = 1.069159849148146
$ws.Range("N6").Value = 1.051150297750944

$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.04500884493704
$ws.Range("D7").Value = 1.058463327839306
$ws.Range("E7").Value = 1.053652546646512
$ws.Range("F7").Value = 1.066568250318308
$ws.Range("I7").Value = 1.049683401567826
$ws.Range("J7").Value = 1.049495694327356
$ws.Range("K7").Value = 1.060891106714365
$ws.Range("L7").Value = 1.056091996224684
$ws.Range("M7").Value = 1.068976622777383
$ws.Range("N7").Value = 1.05098609859214

$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.043859623410829
$ws.Range("D8").Value = 1.057647241302068
$ws.Range("E8").Value = 1.052666365319675
$ws.Range("F8").Value = 1.065559671307905
$ws.Range("I8").Value = 1.049386135989852
$ws.Range("J8").Value = 1.048810325587958
$ws.Range("K8").Value = 1.060319607599443
$ws.Range("L8").Value = 1.055352108707872
$ws.Range("M8").Value = 1.068211065118846
$ws.Range("N8").Value = 1.050299756550519

$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.041839664615783
$ws.Range("D9").Value = 1.056212188099239
$ws.Range("E9").Value = 1.050935618777216
$ws.Range("F9").Value = 1.063788770521911
$ws.Range("I9").Value = 1.048854510968434
$ws.Range("J9").Value = 1.04760264102452
$ws.Range("K9").Value = 1.059310210839916
$ws.Range("L9").Value = 1.054050232232886
$ws.Range("M9").Value = 1.066863285071965
$ws.Range("N9").Value = 1.049090356936477

$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.040496558220531
$ws.Range("D10").Value = 1.055257610827301
$ws.Range("E10").Value = 1.049786623853999
$ws.Range("F10").Value = 1.062612552543502
$ws.Range("I10").Value = 1.048494897491414
$ws.Range("J10").Value = 1.046797609280038
$ws.Range("K10").Value = 1.058635786374438
$ws.Range("L10").Value = 1.053183676965723
$ws.Range("M10").Value = 1.065965686647107
$ws.Range("N10").Value = 1.048284181954579

$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.039915835719327
$ws.Range("D11").Value = 1.054844795749813
$ws.Range("E11").Value = 1.049290266421622
$ws.Range("F11").Value = 1.062104301372361
$ws.Range("I11").Value = 1.048337958982393
$ws.Range("J11").Value = 1.046449059298087
$ws.Range("K11").Value = 1.058343416149489
$ws.Range("L11").Value = 1.052808790032275
$ws.Range("M11").Value = 1.065577255125977
$ws.Range("N11").Value = 1.047935136991675

$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.039700258997045
$ws.Range("D12").Value = 1.054691538350113
$ws.Range("E12").Value = 1.049106074055661
$ws.Range("F12").Value = 1.061915675360019
$ws.Range("I12").Value = 1.04827948206074
$ws.Range("J12").Value = 1.04631959871647
$ws.Range("K12").Value = 1.058234767134025
$ws.Range("L12").Value = 1.052669592439883
$ws.Range("M12").Value = 1.06543301120318
$ws.Range("N12").Value = 1.047805492561183

$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.039746495051935
$ws.Range("D13").Value = 1.054724408909322
$ws.Range("E13").Value = 1.049145575909102
$ws.Range("F13").Value = 1.061956128957603
$ws.Range("I13").Value = 1.048292033825208
$ws.Range("J13").Value = 1.046347368140462
$ws.Range("K13").Value = 1.058258074943337
$ws.Range("L13").Value = 1.052699448414866
$ws.Range("M13").Value = 1.065463950321465
$ws.Range("N13").Value = 1.047833301420942

$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.039898013430039
$ws.Range("D14").Value = 1.054832125792189
$ws.Range("E14").Value = 1.049275037404025
$ws.Range("F14").Value = 1.062088706189319
$ws.Range("I14").Value = 1.04833312898948
$ws.Range("J14").Value = 1.046438357908037
$ws.Range("K14").Value = 1.058334436193775
$ws.Range("L14").Value = 1.052797282844561
$ws.Range("M14").Value = 1.065565331117537
$ws.Range("N14").Value = 1.047924420404424

$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.03999138610627
$ws.Range("D15").Value = 1.054898504442675
$ws.Range("E15").Value = 1.04935482639848
$ws.Range("F15").Value = 1.06217041280491
$ws.Range("I15").Value = 1.048358424856283
$ws.Range("J15").Value = 1.046494420593394
$ws.Range("K15").Value = 1.058381478336509
$ws.Range("L15").Value = 1.052857568813433
$ws.Range("M15").Value = 1.065627800101508
$ws.Range("N15").Value = 1.047980562705225

$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.040535116931559
$ws.Range("D16").Value = 1.055285019189393
$ws.Range("E16").Value = 1.04981959016541
$ws.Range("F16").Value = 1.062646305982305
$ws.Range("I16").Value = 1.048505287277778
$ws.Range("J16").Value = 1.046820742221723
$ws.Range("K16").Value = 1.058655182978373
$ws.Range("L16").Value = 1.05320856420483
$ws.Range("M16").Value = 1.065991470646002
$ws.Range("N16").Value = 1.048307347747695

$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.04087641396079
$ws.Range("D17").Value = 1.055527611090253
$ws.Range("E17").Value = 1.050111437404529
$ws.Range("F17").Value = 1.062945106132344
$ws.Range("I17").Value = 1.048597083225308
$ws.Range("J17").Value = 1.047025445185366
$ws.Range("K17").Value = 1.058826780622605
$ws.Range("L17").Value = 1.053428825642239
$ws.Range("M17").Value = 1.066219655491669
$ws.Range("N17").Value = 1.048512341413027

$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.041075568790023
$ws.Range("D18").Value = 1.0556691611939
$ws.Range("E18").Value = 1.05028177923748
$ws.Range("F18").Value = 1.063119493246102
$ws.Range("I18").Value = 1.048650508141627
$ws.Range("J18").Value = 1.047144848162628
$ws.Range("K18").Value = 1.058926837729437
$ws.Range("L18").Value = 1.053557332851786
$ws.Range("M18").Value = 1.066352774371778
$ws.Range("N18").Value = 1.048631913956211

$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.041143489309454
$ws.Range("D19").Value = 1.055717434604047
$ws.Range("E19").Value = 1.050339880387718
$ws.Range("F19").Value = 1.063178971992357
$ws.Range("I19").Value = 1.048668704598804
$ws.Range("J19").Value = 1.04718556199634
$ws.Range("K19").Value = 1.058960948995386
$ws.Range("L19").Value = 1.053601155947712
$ws.Range("M19").Value = 1.066398168240764
$ws.Range("N19").Value = 1.048672685608237

$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.040839787521899
$ws.Range("D20").Value = 1.055501578053291
$ws.Range("E20").Value = 1.050080113331098
$ws.Range("F20").Value = 1.062913037134969
$ws.Range("I20").Value = 1.048587246599128
$ws.Range("J20").Value = 1.047003482154223
$ws.Range("K20").Value = 1.05880837321196
$ws.Range("L20").Value = 1.053405190308735
$ws.Range("M20").Value = 1.066195171076294
$ws.Range("N20").Value = 1.04849034719186

$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.039853391457646
$ws.Range("D21").Value = 1.054800403617515
$ws.Range("E21").Value = 1.049236909322668
$ws.Range("F21").Value = 1.062049661035683
$ws.Range("I21").Value = 1.048321032531991
$ws.Range("J21").Value = 1.046411563502359
$ws.Range("K21").Value = 1.058311951069538
$ws.Range("L21").Value = 1.052768471593489
$ws.Range("M21").Value = 1.065535475970121
$ws.Range("N21").Value = 1.047897587947616

$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.039233954955167
$ws.Range("D22").Value = 1.054360014109381
$ws.Range("E22").Value = 1.048707777817752
$ws.Range("F22").Value = 1.061507755294181
$ws.Range("I22").Value = 1.0481525949835
$ws.Range("J22").Value = 1.046039438420809
$ws.Range("K22").Value = 1.057999543922409
$ws.Range("L22").Value = 1.052368443452422
$ws.Range("M22").Value = 1.065120913406746
$ws.Range("N22").Value = 1.047524934405767

$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.039562258368102
$ws.Range("D23").Value = 1.054593428050336
$ws.Range("E23").Value = 1.048988182724228
$ws.Range("F23").Value = 1.061794940667439
$ws.Range("I23").Value = 1.048241986937253
$ws.Range("J23").Value = 1.046236704968902
$ws.Range("K23").Value = 1.058165183590282
$ws.Range("L23").Value = 1.052580476813094
$ws.Range("M23").Value = 1.065340660051358
$ws.Range("N23").Value = 1.047722481094985

$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.040856337179734
$ws.Range("D24").Value = 1.05551334110678
$ws.Range("E24").Value = 1.050094266984876
$ws.Range("F24").Value = 1.062927527419454
$ws.Range("I24").Value = 1.048591691712117
$ws.Range("J24").Value = 1.047013406292869
$ws.Range("K24").Value = 1.058816690830308
$ws.Range("L24").Value = 1.05341586999852
$ws.Range("M24").Value = 1.066206234460299
$ws.Range("N24").Value = 1.048500285423922

$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.042361256513359
$ws.Range("D25").Value = 1.056582817329933
$ws.Range("E25").Value = 1.051382213813723
$ws.Range("F25").Value = 1.064245826928296
$ws.Range("I25").Value = 1.048992868482175
$ws.Range("J25").Value = 1.047914845907692
$ws.Range("K25").Value = 1.059571433265874
$ws.Range("L25").Value = 1.054386564601036
$ws.Range("M25").Value = 1.067211562986794
$ws.Range("N25").Value = 1.049403005186393
